$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole row 476 ("「美しい時は作るもの。待つものではない」...")
# which shifts all subsequent rows up by one.
$ws.Rows.Item(476).Delete()
